# Add new power plants to the Electricity Source subscript on the
# "BCRbQ" sheet (rows 19-24), mirroring the formatting of the last
# existing row (row 18) but with literal zero values instead of
# SUMIFS formulas, since Table_6_06 has no data yet for these new
# plant types.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BCRbQ")
$ws.Activate()

$names = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

# Copy the formatting of the last populated row (row 18) down onto
# the six new rows so the new entries match the table's look (font,
# fill, number format, etc.) without carrying over its formulas.
$ws.Range("A18:AF18").Copy()
$ws.Range("A19:AF24").PasteSpecial(-4122)

for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 19 + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]
    for ($c = 2; $c -le 32; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# The header row's explicit height (a leftover autofit value, not a
# user-set custom height) gets recalculated away once the sheet's
# row heights are refreshed after the new rows are added.
$ws.Rows.Item(1).AutoFit()

# Matches the cursor landing one row below the newly entered data,
# as it would after typing the last row by hand.
$ws.Range("A25").Select()
